# Update "想去人数" (want-to-go count) values on both the "展览" sheet
# and the "全部类型" sheet, which mirrors the same data.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F3").Value = 1717
    $ws.Range("F4").Value = 32
    $ws.Range("F8").Value = 79
    $ws.Range("F9").Value = 646
}
